$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case "de/del/la/las/los/el/y" connector words in municipality/state names ---
$ws.Range("B6").Value = "Pabellón De Arteaga"
$ws.Range("B7").Value = "Rincón De Romos"
$ws.Range("B8").Value = "San Francisco De Los Romo"
$ws.Range("B9").Value = "San José De Gracia"
$ws.Range("B31").Value = "Comitán De Domínguez"
$ws.Range("B49").Value = "Salto De Agua"
$ws.Range("B50").Value = "San Cristóbal De Las Casas"
$ws.Range("B73").Value = "Hidalgo Del Parral"
$ws.Range("B83").Value = "Valle De Zaragoza"
$ws.Range("B104").Value = "Villa De Álvarez"
$ws.Range("A106").Value = "Ciudad De México"
$ws.Range("B110").Value = "Cuajimalpa De Morelos"
$ws.Range("B125").Value = "Coneto De Comonfort"
$ws.Range("B137").Value = "Nombre De Dios"
$ws.Range("A150").Value = "Estado De México"
$ws.Range("B150").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B153").Value = "Almoloya De Alquisiras"
$ws.Range("B154").Value = "Almoloya De Juárez"
$ws.Range("B159").Value = "Atizapán De Zaragoza"
$ws.Range("B163").Value = "Chapa De Mota"
$ws.Range("B166").Value = "Coacalco De Berriozábal"
$ws.Range("B171").Value = "Ecatepec De Morelos"
$ws.Range("B174").Value = "Ixtapan De La Sal"
$ws.Range("B175").Value = "Ixtapan Del Oro"
$ws.Range("B183").Value = "Naucalpan De Juárez"
$ws.Range("B190").Value = "San Felipe Del Progreso"
$ws.Range("B193").Value = "Soyaniquilpan De Juárez"
$ws.Range("B201").Value = "Tenango Del Valle"
$ws.Range("B207").Value = "Tlalnepantla De Baz"
$ws.Range("B212").Value = "Valle De Bravo"
$ws.Range("B213").Value = "Villa De Allende"
$ws.Range("B223").Value = "San Miguel De Allende"
$ws.Range("B224").Value = "Apaseo El Alto"
$ws.Range("B225").Value = "Apaseo El Grande"
$ws.Range("B232").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B236").Value = "Jaral Del Progreso"
$ws.Range("B243").Value = "Purísima Del Rincón"
$ws.Range("B247").Value = "San Diego De La Unión"
$ws.Range("B249").Value = "San Francisco Del Rincón"
$ws.Range("B251").Value = "San Luis De La Paz"
$ws.Range("B252").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B254").Value = "Silao De La Victoria"
$ws.Range("B258").Value = "Valle De Santiago"
$ws.Range("B262").Value = "Acapulco De Juárez"
$ws.Range("B264").Value = "Ajuchitlán Del Progreso"
$ws.Range("B265").Value = "Alcozauca De Guerrero"
$ws.Range("B268").Value = "Atenango Del Río"
$ws.Range("B270").Value = "Atoyac De Álvarez"
$ws.Range("B271").Value = "Ayutla De Los Libres"
$ws.Range("B272").Value = "Buenavista De Cuéllar"
$ws.Range("B273").Value = "Chilapa De Álvarez"
$ws.Range("B274").Value = "Chilpancingo De Los Bravo"
$ws.Range("B275").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B280").Value = "Coyuca De Benítez"
$ws.Range("B281").Value = "Coyuca De Catalán"
$ws.Range("B284").Value = "Cuetzala Del Progreso"
$ws.Range("B285").Value = "Cutzamala De Pinzón"
$ws.Range("B291").Value = "Huitzuco De Los Figueroa"
$ws.Range("B292").Value = "Iguala De La Independencia"
$ws.Range("B293").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B294").Value = "Zihuatanejo De Azueta"
$ws.Range("B296").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B307").Value = "Taxco De Alarcón"
$ws.Range("B309").Value = "Técpan De Galeana"
$ws.Range("B311").Value = "Tepecoacuilco De Trujano"
$ws.Range("B313").Value = "Tixtla De Guerrero"
$ws.Range("B324").Value = "Atotonilco El Grande"
$ws.Range("B327").Value = "Cuautepec De Hinojosa"
$ws.Range("B332").Value = "Huejutla De Reyes"
$ws.Range("B335").Value = "Jacala De Ledezma"
$ws.Range("B340").Value = "Mineral Del Chico"
$ws.Range("B341").Value = "Mixquiahuala De Juárez"
$ws.Range("B342").Value = "Molango De Escamilla"
$ws.Range("B344").Value = "Nopala De Villagrán"
$ws.Range("B345").Value = "Pachuca De Soto"
$ws.Range("B347").Value = "Progreso De Obregón"
$ws.Range("B351").Value = "Santiago De Anaya"
$ws.Range("B352").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B356").Value = "Tepehuacán De Guerrero"
$ws.Range("B358").Value = "Tezontepec De Aldama"
$ws.Range("B362").Value = "Tula De Allende"
$ws.Range("B363").Value = "Tulancingo De Bravo"
$ws.Range("B364").Value = "Zacualtipán De Ángeles"
$ws.Range("B365").Value = "Zapotlán De Juárez"
$ws.Range("B368").Value = "Acatlán De Juárez"
$ws.Range("B369").Value = "Ahualulco De Mercado"
$ws.Range("B373").Value = "Atemajac De Brizuela"
$ws.Range("B374").Value = "Atotonilco El Alto"
$ws.Range("B375").Value = "Autlán De Navarro"
$ws.Range("B378").Value = "Cañadas De Obregón"
$ws.Range("B385").Value = "Encarnación De Díaz"
$ws.Range("B390").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B391").Value = "Ixtlahuacán Del Río"
$ws.Range("B400").Value = "La Manzanilla De La Paz"
$ws.Range("B401").Value = "Lagos De Moreno"
$ws.Range("B408").Value = "Ojuelos De Jalisco"
$ws.Range("B413").Value = "San Cristóbal De La Barranca"
$ws.Range("B414").Value = "San Diego De Alejandría"
$ws.Range("B416").Value = "San Juan De Los Lagos"
$ws.Range("B417").Value = "San Juanito De Escobedo"
$ws.Range("B419").Value = "San Martín De Bolaños"
$ws.Range("B421").Value = "San Sebastián Del Oeste"
$ws.Range("B422").Value = "Santa María De Los Ángeles"
$ws.Range("B425").Value = "Talpa De Allende"
$ws.Range("B426").Value = "Tamazula De Gordiano"
$ws.Range("B431").Value = "Teocuitatlán De Corona"
$ws.Range("B432").Value = "Tepatitlán De Morelos"
$ws.Range("B433").Value = "Tlajomulco De Zúñiga"
$ws.Range("B443").Value = "Valle De Guadalupe"
$ws.Range("B444").Value = "Valle De Juárez"
$ws.Range("B448").Value = "Yahualica De González Gallo"
$ws.Range("B449").Value = "Zacoalco De Torres"
$ws.Range("B452").Value = "Zapotlán El Grande"
$ws.Range("B470").Value = "Cojumatlán De Régules"
$ws.Range("B522").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B541").Value = "Coatlán Del Río"
$ws.Range("B548").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B550").Value = "Puente De Ixtla"
$ws.Range("B555").Value = "Tetela Del Volcán"
$ws.Range("B556").Value = "Tlaltizapán De Zapata"
$ws.Range("B565").Value = "Amatlán De Cañas"
$ws.Range("B584").Value = "San Nicolás De Los Garza"
$ws.Range("B587").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B591").Value = "Chalcatongo De Hidalgo"
$ws.Range("B592").Value = "Chiquihuitlán De Benito Juárez"
$ws.Range("B593").Value = "Ciénega De Zimatlán"
$ws.Range("B595").Value = "Coicoyán De Las Flores"
$ws.Range("B596").Value = "Constancia Del Rosario"
$ws.Range("B598").Value = "El Barrio De La Soledad"
$ws.Range("B600").Value = "Guadalupe De Ramírez"
$ws.Range("B601").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B602").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B603").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B604").Value = "Ixtlán De Juárez"
$ws.Range("B605").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B610").Value = "Magdalena Yodocono De Porfirio Díaz"
$ws.Range("B612").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B614").Value = "Oaxaca De Juárez"
$ws.Range("B615").Value = "Ocotlán De Morelos"
$ws.Range("B616").Value = "Pinotepa De Don Luis"
$ws.Range("B617").Value = "Putla Villa De Guerrero"
$ws.Range("B620").Value = "San Agustín De Las Juntas"
$ws.Range("B632").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B664").Value = "San Miguel Del Puerto"
$ws.Range("B675").Value = "San Pablo Villa De Mitla"
$ws.Range("B686").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B704").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B708").Value = "Santa Inés Del Monte"
$ws.Range("B716").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B745").Value = "Santo Domingo De Morelos"
$ws.Range("B753").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B754").Value = "Tanetze De Zaragoza"
$ws.Range("B755").Value = "Tataltepec De Valdés"
$ws.Range("B756").Value = "Teotitlán Del Valle"
$ws.Range("B757").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B758").Value = "Tlacolula De Matamoros"
$ws.Range("B759").Value = "Totontepec Villa De Morelos"
$ws.Range("B761").Value = "Villa De Tututepec"
$ws.Range("B762").Value = "Villa De Zaachila"
$ws.Range("B764").Value = "Villa Sola De Vega"
$ws.Range("B765").Value = "Villa Talea De Castro"
$ws.Range("B767").Value = "Yutanduchi De Guerrero"
$ws.Range("B768").Value = "Zapotitlán Del Río"
$ws.Range("B782").Value = "Chalchicomula De Sesma"
$ws.Range("B793").Value = "Cuayuca De Andrade"
$ws.Range("B801").Value = "Huehuetlán El Chico"
$ws.Range("B806").Value = "Izúcar De Matamoros"
$ws.Range("B813").Value = "Los Reyes De Juárez"
$ws.Range("B817").Value = "Palmar De Bravo"
$ws.Range("B827").Value = "San Nicolás De Los Ranchos"
$ws.Range("B829").Value = "San Salvador El Verde"
$ws.Range("B839").Value = "Tepexi De Rodríguez"
$ws.Range("B841").Value = "Tetela De Ocampo"
$ws.Range("B845").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B856").Value = "Xayacatlán De Bravo"
$ws.Range("B867").Value = "Amealco De Bonfil"
$ws.Range("B869").Value = "Cadereyta De Montes"
$ws.Range("B873").Value = "Landa De Matamoros"
$ws.Range("B877").Value = "San Juan Del Río"
$ws.Range("B885").Value = "Axtla De Terrazas"
$ws.Range("B890").Value = "Ciudad Del Maíz"
$ws.Range("B896").Value = "Mexquitic De Carmona"
$ws.Range("B900").Value = "Santa María Del Río"
$ws.Range("B906").Value = "Tanquián De Escobedo"
$ws.Range("B908").Value = "Villa De Arista"
$ws.Range("B909").Value = "Villa De Guadalupe"
$ws.Range("B910").Value = "Villa De Ramos"
$ws.Range("B911").Value = "Villa De Reyes"
$ws.Range("B947").Value = "Jalpa De Méndez"
$ws.Range("B970").Value = "Soto La Marina"
$ws.Range("B984").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1003").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1006").Value = "Amatlán De Los Reyes"
$ws.Range("B1012").Value = "Boca Del Río"
$ws.Range("B1017").Value = "Castillo De Teayo"
$ws.Range("B1019").Value = "Cazones De Herrera"
$ws.Range("B1031").Value = "Cosamaloapan De Carpio"
$ws.Range("B1045").Value = "Hueyapan De Ocampo"
$ws.Range("B1046").Value = "Ignacio De La Llave"
$ws.Range("B1048").Value = "Ixhuatlán De Madero"
$ws.Range("B1049").Value = "Ixhuatlán Del Café"
$ws.Range("B1050").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1058").Value = "Juchique De Ferrer"
$ws.Range("B1060").Value = "Landero Y Coss"
$ws.Range("B1062").Value = "Lerdo De Tejada"
$ws.Range("B1066").Value = "Martínez De La Torre"
$ws.Range("B1068").Value = "Medellín De Bravo"
$ws.Range("B1071").Value = "Mixtla De Altamirano"
$ws.Range("B1073").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B1084").Value = "Paso De Ovejas"
$ws.Range("B1085").Value = "Paso Del Macho"
$ws.Range("B1088").Value = "Poza Rica De Hidalgo"
$ws.Range("B1094").Value = "Sayula De Alemán"
$ws.Range("B1097").Value = "Soledad De Doblado"
$ws.Range("B1100").Value = "Tatahuicapan De Juárez"
$ws.Range("B1115").Value = "Tlacotepec De Mejía"
$ws.Range("B1124").Value = "Vega De Alatorre"
$ws.Range("B1142").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B1150").Value = "Jiménez Del Teul"
$ws.Range("B1157").Value = "Nochistlán De Mejía"
$ws.Range("B1158").Value = "Noria De Ángeles"
$ws.Range("B1165").Value = "Teúl De González Ortega"
$ws.Range("B1166").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1167").Value = "Trinidad García De La Cadena"

# --- Remove footer / metadata rows (1177-1181) and shrink dimension to A1:D1175 ---
$ws.Range("A1177:A1181").EntireRow.Delete()
